$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3, shifting the existing rows 3-24 down to 4-25
# (dimension grows from A1:R24 to A1:R25).
$ws.Rows(3).Insert()

# Populate the newly inserted row 3 with the new weekly price entry.
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 44530
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 100112028
$ws.Range("G3").Value = "Sandia"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 1300
$ws.Range("K3").Value = 450
$ws.Range("L3").Value = 480
$ws.Range("M3").Value = 465
$ws.Range("N3").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O3").Value = "Perú"
$ws.Range("P3").Value = 465
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = "Hortaliza"
